# Update data + results
# update data (2020) + results new mod (GB reviewed)
#
# The "Clim_year" and "Clim_season" worksheets swap names (the tab that used
# to be called "Clim_year" is now "Clim_season" and vice versa, while the
# underlying sheet order / data stay where they are), and the previously
# active "Summary" tab is replaced by the (renamed) second tab as the
# selected/active sheet.

$wb = $excel.ActiveWorkbook

# Swap the names of the "Clim_year" and "Clim_season" sheets via a temporary
# name so neither ever collides with the other.
$wsYear = $wb.Worksheets.Item("Clim_year")
$wsSeason = $wb.Worksheets.Item("Clim_season")

$wsYear.Name = "Clim_year_tmp"
$wsSeason.Name = "Clim_year"
$wsYear.Name = "Clim_season"

# The sheet now named "Clim_season" (previously "Clim_year", second tab)
# becomes the active/selected sheet instead of "Summary".
$wb.Worksheets.Item("Clim_season").Activate()
